# Add the three new character styles that the commit introduces, then
# apply them to the runs in the newly-added paragraphs.

$d = $word.ActiveDocument
$styles = $d.Styles

# --- GaNStyle: Calibri, 14pt --------------------------------------------
$gaNStyle = $styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

# --- GaNParagraph: Calibri, 10pt ----------------------------------------
$gaNParagraph = $styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

# --- GaNLinks: Calibri, 9.5pt, bold, navy, underline --------------------
$gaNLinks = $styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Underline = 1

# Apply GaNStyle to every "2022 Campaign Dates..." run (4 occurrences).
$rng = $d.Content
while ($rng.Find.Execute(" 2022 Campaign Dates that use Taurus constellation: January 16-25", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
  $rng.Style = "GaNStyle"
  $rng.Collapse(0)
  $rng.SetRange($rng.End, $d.Content.End)
}

# Apply GaNParagraph to the campaign description paragraph.
$rng2 = $d.Content
if ($rng2.Find.Execute("You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Taurus constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
  $rng2.Style = "GaNParagraph"
}

# Apply GaNLinks to the map URL run.
$rng3 = $d.Content
if ($rng3.Find.Execute("(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
  $rng3.Style = "GaNLinks"
}

Write-Output "done"
